# Send tutoring request code
#
# Adds two new Product Backlog rows (20 and 21):
#   <SP19> / Comments / Teacher can filter tutors by subject in comment tutor page
#   <SP20> / Emails for admin / Admin can send E-Mails to an outstanding tutor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same formatting used by the existing data rows (row 19) to the
# two new rows before filling in values.
$ws.Range("A19:I19").Copy() | Out-Null
$ws.Range("A20:I20").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:I21").PasteSpecial(-4122) | Out-Null

# Match the row height used throughout the sheet.
$ws.Rows.Item(20).RowHeight = 15.75
$ws.Rows.Item(21).RowHeight = 15.75

# Fill in the values column by column (A column first for both rows, then B,
# then C) so new shared-string entries are created in the same order Excel
# produced them in.
$ws.Range("A20").Value = "<SP19>"
$ws.Range("A21").Value = "<SP20>"

$ws.Range("B20").Value = "Comments"
$ws.Range("C20").Value = "Teacher can filter tutors by subject in comment tutor page"

$ws.Range("B21").Value = "Emails for admin"
$ws.Range("C21").Value = "Admin can send E-Mails to an outstanding tutor"

$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 1
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 3

$ws.Range("F20").Value = "no"
$ws.Range("G20").Value = "no"
$ws.Range("H20").Value = "no"
$ws.Range("I20").Value = "no"

$ws.Range("F21").Value = "no"
$ws.Range("G21").Value = "no"
$ws.Range("H21").Value = "no"
$ws.Range("I21").Value = "no"

# Stray formatted-but-empty cells J20:Q20 (artifact of an over-wide paste),
# matching the committed workbook.
$ws.Range("A17").Copy() | Out-Null
$ws.Range("J20:Q20").PasteSpecial(-4122) | Out-Null
$ws.Range("J20:Q20").Value = $null
$ws.Range("J20:Q20").Borders.LineStyle = -4142

# Selection as recorded in the edited workbook.
$ws.Range("C11").Select() | Out-Null
